$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$row = 24
$ws.Cells.Item($row, 1).Value = 2
$ws.Cells.Item($row, 2).Value = "Organized crime-wave"
$ws.Cells.Item($row, 3).Value = "Nick Fury"
$ws.Cells.Item($row, 4).Value = "Avengers|Illuminati"
$ws.Cells.Item($row, 5).Value = "Maggia Goons"
$ws.Cells.Item($row, 6).Value = "Speed (R)|The Captain & The Devil (SW2)|Goliath (CW)|Captain Marvel, Agent of SHIELD (R)|Totally Awesome Hulk (CH)"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "2|9"
$ws.Cells.Item($row, 9).Value = "no"
$ws.Cells.Item($row, 11).Value = "Goons sequenced into a twist once for 3 and once for 2. Lots of wounds and insufficient thinning as goons get twisted before you can beat them."

$ws.Range("K25").Select()
